$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-15 from 2023-10-22 to 2023-10-25
$ws.Range("C2:C15").Value = 45224
